$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row gets a thin box border and top-vertical alignment (new source formatting)
$ws.Range("A1:D1").Borders.LineStyle = 1
$ws.Range("A1:D1").VerticalAlignment = -4160

# Update A column (region names) - reorder rows
$ws.Range("A2").Value = "Tocantins"
$ws.Range("A3").Value = "Maranhão"
$ws.Range("A4").Value = "Piauí"
$ws.Range("A5").Value = "Alagoas"
$ws.Range("A6").Value = "Amazonas"
$ws.Range("A7").Value = "Rondônia"
$ws.Range("A8").Value = "Sergipe"
$ws.Range("A9").Value = "Brasil"
$ws.Range("A10").Value = "Nordeste"

# Update B column (variable label text) for all data rows
$ws.Range("B2:B10").Value = "Diferença 2010-2000"

# Update C column (values) to match new figures
$ws.Range("C2").Value = 0.1739999999999999
$ws.Range("C3").Value = 0.163
$ws.Range("C4").Value = 0.162
$ws.Range("C5").Value = 0.16
$ws.Range("C6").Value = 0.159
$ws.Range("C7").Value = 0.1529999999999999
$ws.Range("C8").Value = 0.147
$ws.Range("C9").Value = 0.115
$ws.Range("C10").Value = 0.119

# Update D column (ranking) - Sergipe ranking changed 10º -> 9º
# (D9 and D10 have no ranking, same as before - Nordeste/Brasil rows)
$ws.Range("D8").Value = "9º"
